# Actualización automática del inventario, Google Sheets y productos.json
# Adds a new inventory row (row 33) for "Correa de transporte Epson".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33

$ws.Cells.Item($row, 1).Value = "CJVMIV"
$ws.Cells.Item($row, 2).Value = "Correa de transporte Epson"
$ws.Cells.Item($row, 3).Value = "TM U950"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 7
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Formula = "=(E33-D33)*G33"
$ws.Cells.Item($row, 9).Formula = "=D33*F33"
$ws.Cells.Item($row, 10).Value = 0
